$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 392, shifting existing rows 392:479 down to 393:480.
$ws.Rows.Item(392).Insert()

# Populate the newly inserted row 392 with its data (matching the other rows'
# constant columns, plus the new record's own values).
$ws.Range("A392").Value = 4
$ws.Range("B392").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C392").Value = "Los Lagos"
$ws.Range("D392").Value = 45275
$ws.Range("E392").Value = 10
$ws.Range("F392").Value = 100112044
$ws.Range("G392").Value = "Perejil"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 120
$ws.Range("K392").Value = 6000
$ws.Range("L392").Value = 6000
$ws.Range("M392").Value = 6000
$ws.Range("N392").Value = "$/docena de atados (2 kilos)"
$ws.Range("O392").Value = "Región de La Araucanía"
$ws.Range("P392").Value = 3000
$ws.Range("Q392").Value = 2
$ws.Range("R392").Value = "Hortaliza"
